# Saldo.xlsx update: re-order the top balances, update EVANGELINA's saldo,
# and drop several accounts that are no longer present in the export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete rows that disappear entirely (bottom-to-top so row numbers
#        of earlier rows stay valid while we work). These are, in the
#        original layout: LAGO, GUILHERME, EDMAR, SIMONE (rows 3-6, which
#        get re-inserted above EVANGELINA with the same values), MVFC (row 8,
#        removed outright), the old MONICA row (row 10, replaced by a new
#        value above), MERG (row 17, removed outright), and JOSE/1097.39 +
#        BRUNO/981.31 (rows 22-23, removed outright).
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# --- 2. Update EVANGELINA's balance (row 2, unchanged position).
$ws.Cells.Item(2, 3).Value = 138881.35

# --- 3. Insert 5 fresh rows above EVANGELINA (row 2) and populate them.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(2).Insert()
}

$newRows = @(
    @("001882235", "LAGO", 289624.21),
    @("005142611", "GUILHERME", 174663.16),
    @("008004851", "EDMAR", 150010.04),
    @("008004870", "SIMONE", 150010.04),
    @("004387250", "MONICA", 143045.29)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 2 + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}
